$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row for the "102_AutomobileInsurance_002_VehicleData_002_FieldHintsAndErrors"
# record, right after the existing VehicleData_001 row (row 7), then re-sort the data
# range by the "Record/Process" column so everything lines back up in order.
$ws.Rows("4:4").Insert() | Out-Null

$ws.Range("A4").Value = "102_AutomobileInsurance_002_VehicleData_002_FieldHintsAndErrors"
$ws.Range("B4").Value = "var102_AutomobileInsurance_002_VehicleData_002_FieldHintsAndErrors"
$ws.Range("C4").Value = "Open Automobile Insurance"
$ws.Range("D4").Value = "102_AutomobileInsurance_002_VehicleData_002_FieldHintsAndErrors"

# Re-sort rows 2:9 (the data, excluding the header row) ascending by column A so the
# newly inserted record lands in its natural alphabetical/numeric position.
$dataRange = $ws.Range("A2:G9")
$keyRange = $ws.Range("A2:A9")
$dataRange.Sort($keyRange)

# Widen columns A, B and D to fit the new, longer record/variable names.
$ws.Columns("A").ColumnWidth = 71.3
$ws.Columns("B").ColumnWidth = 60.6
$ws.Columns("D").ColumnWidth = 60.6

# Restore the last-used selection.
$ws.Range("B18").Select() | Out-Null
